$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell format (bold, border, centered) from B1 to A2
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Set the values for the new row
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0.5864994333333335
$ws.Range("C2").Value = -1.147916233333333
$ws.Range("D2").Value = -0.2273942333333316
